$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove extra columns U:AD (rows 1-19) - shrinks the grid width from AD to T
$ws.Range("U1:AD19").Delete() | Out-Null

# 2. Update row 2 header cells C2:T2 with the new shared-string order
$row2Headers = @('[2, 1, 1]', '[4, 0, 0]', '[2, 0, 0]', '[2, 2, 0]', '[1, 1, 0]', '[3, 1, 0]', '[2, 2, 2]', '[3, 2, 1]', '1Pair-A', '1Pair-B', '2Pairs-A', '2Pairs-B', '3Pairs-A', '3Pairs-B', '3Pairs-C', '4Pairs', '5A4F', 'MaxUnique')
for ($i = 0; $i -lt $row2Headers.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $row2Headers[$i]
}

# 3. Row labels (column B) for rows 3-29; $labels[0] is row 3 .. $labels[26] is row 29
$labels = @('Spiral5', 'RotRing OmegaMax-90', 'Equal Angle', 'Tilt Rotate', 'CLR', 'Rizzie Hex', 'Thomas Hex', 'Tilt Rotate_Partial', 'RotRing OmegaMax-60', 'Equal Angle_Partial', 'Rizzie Hex_Partial', 'ND Single', 'RD Single', 'TD Single', 'Morris Single', 'Ring Perpendicular to ND', 'Ring Perpendicular to RD', 'Ring Perpendicular to TD', 'OffsetFTD', 'OffsetATD', 'OffsetF45', 'OffsetA45', 'OffsetFRD', 'OffsetARD', 'Gaussian Quadrature', 'Michael-CCHex', 'Michael-SNHex')
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 3
    $ws.Cells.Item($r, 2).Value = $labels[$i]
}

# 4. Add new rows 20-29: column A index, column B label (from $labels), C:T = 1
#    Formatting for column A is copied from the existing A19 cell (bold/border/center style).
for ($r = 20; $r -le 29; $r++) {
    $aValue = $r - 2
    $ws.Cells.Item(19, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $aValue
    $ws.Cells.Item($r, 2).Value = $labels[$aValue - 1]
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
$excel.CutCopyMode = 0
